$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = 6557
$ws.Range("C27").Value = 1016
$ws.Range("D27").Value = 6118000
$ws.Range("E27").Value = 933.0486502973921
$ws.Range("F27").Value = 10.2016806722689
$ws.Range("G27").Value = 7.286166842661035
$ws.Range("H27").Value = 25.48438552134167
